$wb = $excel.ActiveWorkbook

# Add violent crime data for week ending 2023-05-04 (column J = year-to-date 2023 totals)
# 189 cell updates across 47 worksheets, generated from the source diff.

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 2306
$ws.Range("H3").Value = 8347
$ws.Range("I3").Value = 7486
$ws.Range("J3").Value = 2395
$ws.Range("B4").Value = 1672
$ws.Range("D4").Value = 1952
$ws.Range("E4").Value = 1987
$ws.Range("J4").Value = 544
$ws.Range("J6").Value = 3014
$ws.Range("B7").Value = 23304
$ws.Range("D7").Value = 28142
$ws.Range("E7").Value = 25991
$ws.Range("H7").Value = 26005
$ws.Range("I7").Value = 26203
$ws.Range("J7").Value = 8428

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 65
$ws.Range("J4").Value = 35
$ws.Range("J7").Value = 259
$ws.Range("J8").Value = 533
$ws.Range("J10").Value = 51
$ws.Range("J11").Value = 115
$ws.Range("J14").Value = 30
$ws.Range("J15").Value = 103
$ws.Range("J18").Value = 98
$ws.Range("J19").Value = 269
$ws.Range("J20").Value = 178
$ws.Range("H23").Value = 262
$ws.Range("J23").Value = 78
$ws.Range("J25").Value = 50
$ws.Range("J26").Value = 12
$ws.Range("J29").Value = 470
$ws.Range("J30").Value = 34
$ws.Range("J33").Value = 349
$ws.Range("J36").Value = 130
$ws.Range("J37").Value = 287
$ws.Range("J42").Value = 319
$ws.Range("J43").Value = 80
$ws.Range("J47").Value = 75
$ws.Range("J48").Value = 80
$ws.Range("J50").Value = 46
$ws.Range("J51").Value = 112
$ws.Range("J52").Value = 205
$ws.Range("J54").Value = 170
$ws.Range("B63").Value = 376
$ws.Range("D63").Value = 334
$ws.Range("E63").Value = 332
$ws.Range("J63").Value = 41
$ws.Range("I65").Value = 611
$ws.Range("J65").Value = 221
$ws.Range("J67").Value = 308
$ws.Range("J70").Value = 16
$ws.Range("J76").Value = 119
$ws.Range("J77").Value = 62
$ws.Range("J78").Value = 116
$ws.Range("J79").Value = 257
$ws.Range("J83").Value = 198
$ws.Range("J84").Value = 81
$ws.Range("J85").Value = 398
$ws.Range("J89").Value = 88
$ws.Range("J90").Value = 93
$ws.Range("J92").Value = 27
$ws.Range("J94").Value = 69
$ws.Range("J95").Value = 125
$ws.Range("J96").Value = 95
$ws.Range("J99").Value = 116
$ws.Range("B101").Value = 23304
$ws.Range("D101").Value = 28142
$ws.Range("E101").Value = 25991
$ws.Range("H101").Value = 26005
$ws.Range("I101").Value = 26203
$ws.Range("J101").Value = 8428

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 97
$ws.Range("J6").Value = 116
$ws.Range("J7").Value = 398

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J7").Value = 205

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 39
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 168
$ws.Range("J3").Value = 172
$ws.Range("J6").Value = 155
$ws.Range("J7").Value = 533

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 86
$ws.Range("J3").Value = 77
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 30
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J3").Value = 6
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 84
$ws.Range("J3").Value = 105
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 287

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 36
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 124
$ws.Range("J7").Value = 308

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J3").Value = 20
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 182
$ws.Range("J6").Value = 80
$ws.Range("I7").Value = 611
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J3").Value = 69
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 45
$ws.Range("J3").Value = 36
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 125

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 102
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 349

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 33
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 137
$ws.Range("J3").Value = 163
$ws.Range("J7").Value = 470

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 66
$ws.Range("J6").Value = 101
$ws.Range("J7").Value = 269

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 16
$ws.Range("J4").Value = 15
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 119

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J6").Value = 165
$ws.Range("J7").Value = 319

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value = 40
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 22
$ws.Range("H3").Value = 80
$ws.Range("J6").Value = 21
$ws.Range("H7").Value = 262
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 96
$ws.Range("J7").Value = 257

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 55
$ws.Range("J3").Value = 53
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 48
$ws.Range("J3").Value = 32
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value = 14
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 21
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 12

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 24
$ws.Range("J6").Value = 19

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J3").Value = 20
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 112

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 35
